$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.640.60'
$ws.Range("E2").Value = '  +5.56%  '

$ws.Range("D3").Value = '2.626.71'
$ws.Range("E3").Value = '  +10.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.71'
$ws.Range("E5").Value = '  +6.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.99'
$ws.Range("E6").Value = '  +12.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("E7").Value = '  +9.77%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").Value = '  +19.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.13'
$ws.Range("E10").Value = '  +17.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '55.42'
$ws.Range("E11").Value = '  +3.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0859'
$ws.Range("E12").Value = '  +10.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.34'
$ws.Range("E13").Value = '  +19.76%  '

$ws.Range("D14").Value = '3.036.04'
$ws.Range("E14").Value = '  +10.80%  '

$ws.Range("E15").Value = '  +3.52%  '

$ws.Range("D16").Value = '2.635.35'
$ws.Range("E16").Value = '  +10.69%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.941'
$ws.Range("E17").Value = '  +14.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.30'
$ws.Range("E18").Value = '  +9.72%  '

$ws.Range("D19").Value = '47.701.05'
$ws.Range("E19").Value = '  +5.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000103'
$ws.Range("E20").Value = '  +11.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.16'
$ws.Range("E21").Value = '  +6.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.76'
$ws.Range("E22").Value = '  +11.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.11'
$ws.Range("E23").Value = '  +10.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '272.17'
$ws.Range("E24").Value = '  +14.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.10'
$ws.Range("E25").Value = '  +12.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.24'
$ws.Range("E26").Value = '  +19.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '30.48'
$ws.Range("E27").Value = '  +45.85%  '

$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("E29").Value = '  +1.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.72'
$ws.Range("E30").Value = '  +12.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '40.06'
$ws.Range("E31").Value = '  +7.34%  '

$ws.Range("E32").Value = '  +4.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.17'
$ws.Range("E33").Value = '  +14.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.68'
$ws.Range("E34").Value = '  -2.76%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.89'
$ws.Range("E35").Value = '  +7.31%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0857'
$ws.Range("E36").Value = '  +13.41%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.22'
$ws.Range("E37").Value = '  +14.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.23'
$ws.Range("E38").Value = '  +3.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.125'
$ws.Range("E39").Value = '  +12.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.124'
$ws.Range("E40").Value = '  +9.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.45'
$ws.Range("E41").Value = '  +14.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.03'
$ws.Range("E42").Value = '  +55.58%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.28'
$ws.Range("E43").Value = '  +16.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.69'
$ws.Range("E44").Value = '  +17.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0334'
$ws.Range("E45").Value = '  +14.39%  '

$ws.Range("D46").Value = '2.201.69'
$ws.Range("E46").Value = '  +11.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '95.90'
$ws.Range("E47").Value = '  +8.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  +0.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.13'
$ws.Range("E49").Value = '  +20.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '114.58'
$ws.Range("E50").Value = '  +15.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.80'
$ws.Range("E51").Value = '  +6.34%  '
